$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '24.606.82'
$ws.Range('E2').Value = '  +1.64%  '
$ws.Range('D3').Value = '1.704.28'
$ws.Range('E3').Value = '  +1.65%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.004'
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '308.76'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9997'
$ws.Range('E6').Value = '  +0.08%  '
$ws.Range('E7').Value = '  -0.14%  '
$ws.Range('E8').Value = '  +2.68%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3416'
$ws.Range('E9').Value = '  -1.03%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.178'
$ws.Range('E10').Value = '  -0.63%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07420'
$ws.Range('E11').Value = '  +1.48%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.001'
$ws.Range('E12').Value = '  +0.08%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '20.76'
$ws.Range('E13').Value = '  +1.54%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.192'
$ws.Range('E14').Value = '  +1.27%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.892'
$ws.Range('E15').Value = '  +1.62%  '
$ws.Range('D16').Value = '1.702.01'
$ws.Range('E16').Value = '  +1.42%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001114'
$ws.Range('E17').Value = '  +0.30%  '
$ws.Range('E18').Value = '  +0.13%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06659'
$ws.Range('E19').Value = '  -1.02%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '82.94'
$ws.Range('E20').Value = '  +1.29%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.01'
$ws.Range('E21').Value = '  +2.59%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.317'
$ws.Range('E22').Value = '  +2.65%  '
$ws.Range('E23').Value = '  +9.19%  '
$ws.Range('D24').Value = '24.564.82'
$ws.Range('E24').Value = '  +1.60%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.436'
$ws.Range('E25').Value = '  +1.05%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.754'
$ws.Range('E26').Value = '  +3.26%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.06'
$ws.Range('E27').Value = '  +2.54%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '148.88'
$ws.Range('E28').Value = '  -1.73%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '130.57'
$ws.Range('E29').Value = '  +2.96%  '
$ws.Range('D30').Value = '1.886.91'
$ws.Range('E30').Value = '  +1.29%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.166'
$ws.Range('E31').Value = '  +17.23%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.650'
$ws.Range('E32').Value = '  +3.15%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.210'
$ws.Range('E33').Value = '  +2.89%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.08749'
$ws.Range('E34').Value = '  +2.77%  '
$ws.Range('B35').Value = 'WEMIXTOKEN'
$ws.Range('C35').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.756'
$ws.Range('E35').Value = '  -1.14%  '
$ws.Range('B36').Value = 'Aptos'
$ws.Range('C36').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '13.49'
$ws.Range('E36').Value = '  +6.57%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.476'
$ws.Range('E37').Value = '  +1.84%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.06474'
$ws.Range('E38').Value = '  -0.29%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '8.861'
$ws.Range('E39').Value = '  -0.54%  '
$ws.Range('E40').Value = '  +0.17%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.2176'
$ws.Range('E41').Value = '  +1.38%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.273'
$ws.Range('E42').Value = '  -0.39%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.6367'
$ws.Range('E43').Value = '  +2.73%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.9995'
$ws.Range('E44').Value = '  +0.08%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.83'
$ws.Range('E45').Value = '  +4.17%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.6042'
$ws.Range('E46').Value = '  +1.51%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.794'
$ws.Range('E47').Value = '  -0.45%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.093'
$ws.Range('E48').Value = '  +2.92%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '128.24'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.07220'
$ws.Range('E50').Value = '  +0.68%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '78.54'
$ws.Range('E51').Value = '  +2.47%  '
